$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) price cells that look numeric to remain text,
# matching the workbook author style (inline/shared text strings, not numbers).
$ws.Range("D2:D27").NumberFormat = "@"
$ws.Range("D40:D45").NumberFormat = "@"
$ws.Range("D47:D50").NumberFormat = "@"

$ws.Range("D2").Value = '249.82'
$ws.Range("E2").Value = '1BNBBNBBestin24h'
$ws.Range("D3").Value = '22.32'
$ws.Range("D4").Value = '5.630'
$ws.Range("D5").Value = '0.05599'
$ws.Range("D6").Value = '3.373'
$ws.Range("D7").Value = '6.482'
$ws.Range("D8").Value = '1.082'
$ws.Range("D9").Value = '0.8003'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1421'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.07461'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = '0.03285'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'ProBitToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D13").Value = '0.1290'
$ws.Range("E13").Value = '12ProBitTokenPROB'
$ws.Range("D14").Value = '0.02990'
$ws.Range("D15").Value = '0.09261'
$ws.Range("D16").Value = '0.001662'
$ws.Range("D17").Value = '3.246'
$ws.Range("D18").Value = '0.04727'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").Value = '0.0005728'
$ws.Range("E19").Value = '18OneONEWorstin24h'
$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D20").Value = '0.006246'
$ws.Range("E20").Value = '19TigerCashTCH'
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").Value = '0.001053'
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("B22").Value = 'HotbitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D22").Value = '0.003825'
$ws.Range("E22").Value = '21HotbitTokenHTB'
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").Value = '0.0001497'
$ws.Range("E23").Value = '22NitroExNTX'
$ws.Range("B24").Value = 'UpBots'
$ws.Range("C24").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D24").Value = '0.0004766'
$ws.Range("E24").Value = '23UpBotsUBXT'
$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D25").Value = '3.983'
$ws.Range("E25").Value = '24LEOLEO'
$ws.Range("B26").Value = 'BTSEToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D26").Value = '2.136'
$ws.Range("E26").Value = '25BTSETokenBTSE'
$ws.Range("B27").Value = 'BitpandaEcosystemToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D27").Value = '0.3311'
$ws.Range("E27").Value = '26BitpandaEcosystemTokenBEST'
$ws.Range("D40").Value = '0.04209'
$ws.Range("D41").Value = '0.007045'
$ws.Range("D42").Value = '0.1048'
$ws.Range("D43").Value = '0.003085'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = '0.009016'
$ws.Range("D45").Value = '0.00005625'
$ws.Range("D47").Value = '0.6784'
$ws.Range("D48").Value = '0.03079'
$ws.Range("D49").Value = '0.00002096'
$ws.Range("D50").Value = '0.01008'
